$d = $word.ActiveDocument

# The document contains four paragraphs with the (Spanish) "campaign dates"
# text for "Perseo" built up out of many small runs (plus interleaved
# w:proofErr markers from Word's grammar/spell checker). Each of those
# paragraphs must be collapsed down to a single run containing the new,
# already-translated sentence - with no leftover rPr/proofErr artifacts.

$newText = "Fechas de la campaña para Perseo: Hercules: 13-22 de junio, 12-21 de julio, del 10 al 19 de agosto"

function Replace-ParagraphText($paragraph, $text) {
    $full = $paragraph.Range
    # Exclude the paragraph mark itself so we never merge this paragraph
    # into the next one.
    $full.MoveEnd(1, -1)
    $contentEnd = $full.End

    # Some of these paragraphs end with a trailing w:proofErr element that
    # sits right before the paragraph mark. A plain Range.Delete() on the
    # exact content range leaves that trailing marker behind (it only gets
    # cleaned up when the deletion also consumes a run that follows it).
    # Work around this by appending a harmless one-character placeholder
    # run right after the real content (this pushes any trailing proofErr
    # marker off of the paragraph-end boundary), then delete the original
    # content *together with* that placeholder in one shot - which reliably
    # sweeps up every proofErr marker in between.
    $placeholderPoint = $d.Range($contentEnd, $contentEnd)
    $placeholderPoint.InsertAfter("X")

    $deleteRange = $d.Range($full.Start, $contentEnd + 1)
    $deleteRange.Delete()

    # Paragraph is now completely empty (no runs, no proofErr). Insert the
    # replacement text; on an empty range this creates a single new run
    # with no rPr at all, matching a from-scratch <w:r><w:t>...</w:t></w:r>.
    $insertion = $paragraph.Range
    $insertion.MoveEnd(1, -1)
    $insertion.InsertAfter($text)
}

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($t -like "*Fechas*" -and $t -like "*Perseo*" -and $t -like "*diciembre*") {
        Replace-ParagraphText $p $newText
    }
}

Write-Output "done"
